$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column J (10) for the new "서비스 ID" search field.
# This shifts old J/K/L (액티비티 ID / Exception Text / Exception Stack) to K/L/M,
# and the first column-width band grows from 1-10 to 1-11 as a result.
$ws.Columns.Item(10).Insert()

# New label cell (row 4) + its input cell, copying the formatting of the
# existing label/input pairs (G4 label style, H4 input style).
$ws.Range("G4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("I4").Value = "에러코드"

# Row 5's first label changes from "에러코드" to "서비스 ID", and a new blank
# input cell appears under the new I4 label (style copied from H5).
$ws.Range("A5").Value = "서비스 ID"
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)

# New header cell for the inserted column in the results header row (row 7).
$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = "서비스 ID"

$excel.CutCopyMode = $false

# Selection moves to the newly added input range.
$ws.Range("I4:J4").Select()
